$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table runs through row 1106 (A1:R1106). A new record needs to be
# inserted as row 1057, pushing the existing rows 1057-1106 down to
# 1058-1107 (dimension becomes A1:R1107).
$ws.Rows.Item(1057).Insert()

# Populate the newly inserted row with the new "Cuatro cascos verde" record.
$ws.Cells.Item(1057, 1).Value = 5
$ws.Cells.Item(1057, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1057, 3).Value = "Maule"
$ws.Cells.Item(1057, 4).Value = 45267
$ws.Cells.Item(1057, 5).Value = 7
$ws.Cells.Item(1057, 6).Value = 100112002
$ws.Cells.Item(1057, 7).Value = "Pimiento"
$ws.Cells.Item(1057, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(1057, 9).Value = "Primera"
$ws.Cells.Item(1057, 10).Value = 300
$ws.Cells.Item(1057, 11).Value = 14000
$ws.Cells.Item(1057, 12).Value = 14000
$ws.Cells.Item(1057, 13).Value = 14000
$ws.Cells.Item(1057, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(1057, 15).Value = "Región del Maule"
$ws.Cells.Item(1057, 16).Value = 933
$ws.Cells.Item(1057, 17).Value = 15
$ws.Cells.Item(1057, 18).Value = "Hortaliza"
